$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table "Condicion_Pacientes" currently spans A1:F87 (86 data rows).
# Add one new data row (88) for 06/08/2020 and let the table auto-expand.
$tbl = $ws.ListObjects.Item(1)
$newRow = $tbl.ListRows.Add()

# Fill in the new row's values (Fecha, Pruebas Realizadas, Pruebas Positivas,
# Clinicamente Estables, Clinicamente Graves, Cuidados Intensivos).
$ws.Cells.Item(88, 1).Value = 43990
$ws.Cells.Item(88, 2).Value = 682
$ws.Cells.Item(88, 3).Value = 123
$ws.Cells.Item(88, 4).Value = 389
$ws.Cells.Item(88, 5).Value = 215
$ws.Cells.Item(88, 6).Value = 53

# Match the formatting of the row above (date format in column A, centered
# number format in B:F) by copying formats down from row 87.
$ws.Range("A87:F87").Copy()
$ws.Range("A88:F88").PasteSpecial(-4122)

# Move the selection to the new last cell, matching the saved selection.
$ws.Range("F88").Select()
